$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---- Shape 1: "Heading" textbox ----
$sh1 = $s.Shapes.Item(1)
$sh1.Left = 15.75
$sh1.Width = 704.25
$sh1.Height = 21.6

$tr1 = $sh1.TextFrame.TextRange
$tr1.ParagraphFormat.SpaceBefore = 0
$tr1.ParagraphFormat.SpaceAfter = 15
$tr1.Font.Size = 12
$tr1.Font.Bold = $true
$tr1.Font.Name = "Arial"

# ---- Shape 2: "Very long line ..." textbox ----
$sh2 = $s.Shapes.Item(2)
$sh2.Left = 15.75
$sh2.Top = 66
$sh2.Width = 704.25
$sh2.Height = 73.49992125984252

$tr2 = $sh2.TextFrame.TextRange
$tr2.ParagraphFormat.SpaceBefore = 0
$tr2.ParagraphFormat.SpaceAfter = 0
$tr2.Font.Size = 10
$tr2.Font.Name = "Arial"
